$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50:56 down to 51:57
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new weekly price record
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44522
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 100112040
$ws.Range("G50").Value = "Cilantro"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 300
$ws.Range("K50").Value = 1400
$ws.Range("L50").Value = 1500
$ws.Range("M50").Value = 1450
$ws.Range("N50").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 725
$ws.Range("Q50").Value = 2
$ws.Range("R50").Value = "Hortaliza"
